$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (new agenda entry: Pedro / Lindolfo Mourão, concluded) ---
$ws.Range("A12").Value = "Pedro"
$ws.Range("B12").Value = "'2685"
$ws.Range("C12").Value = "Lindolfo Mourão"
$ws.Range("D12").Value = "Remanejamento de um sensor."
$ws.Range("E12").Value = "Técnico esteve no local e disse que o cliente não tem interesse em mexer nisso agora."
$ws.Range("G12").Value = "Concluido"

# --- Row 13 (Pedro / Pirobras, pending) ---
$ws.Range("A13").Value = "Pedro"
$ws.Range("B13").Value = "'2503"
$ws.Range("C13").Value = "Pirobras"
$ws.Range("D13").Value = "Cliente pedindo reparo em zonas com disparo frequente."
$ws.Range("G13").Value = "Pendente"

# --- Row 14 (Pedro / Marco Otávio, pending; long text wraps -> taller row) ---
$ws.Range("A14").Value = "Pedro"
$ws.Range("B14").Value = "'2029"
$ws.Range("C14").Value = "Marco Otávio"
$ws.Range("D14").WrapText = $true
$ws.Range("D14").Value = "Cliente pedindo revisão no alarme, ele acha que está com problema (disse que pode ser na bateria)."
$ws.Rows.Item(14).RowHeight = 30
$ws.Range("G14").Value = "Pendente"

# --- Row 15 (Pedro / Anselmo, pending) ---
$ws.Range("A15").Value = "Pedro"
$ws.Range("B15").Value = "'2194"
$ws.Range("C15").Value = "Anselmo"
$ws.Range("D15").Value = "Cliente pedindo pra ver sobre a instabilidade no alarme dele (é gprs)."
$ws.Range("G15").Value = "Pendente"

# --- Row 16 (Pedro / Casa da Luci, concluded/removed) ---
$ws.Range("A16").Value = "Pedro"
$ws.Range("B16").Value = "'2361"
$ws.Range("C16").Value = "Casa da Luci"
$ws.Range("D16").Value = "Sem comunicação de alarmes."
$ws.Range("E16").Value = "Foi removido completamente do sistema, cliente disse que não tem nada conosco."
$ws.Range("G16").Value = "Concluido"

# Update the saved selection / scroll state to the last-touched cell,
# matching the author's final cursor position (also clears the stale
# topLeftCell freeze left over from the previous save).
$ws.Range("H16").Select()
